$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11 gets the label that used to belong to row 10 ("2020 (14.01.21)"),
# written first so the shared-string table keeps this text at the slot the
# renamed string previously occupied.
$ws.Range("A11").Value = "2020 (14.01.21)"
$ws.Range("B11").Value = 335
$ws.Range("C11").Value = 62
$ws.Range("D11").Value = 23
$ws.Range("E11").Value = 146
$ws.Range("F11").Value = 160

# Row 10 keeps its original numeric data but gets a new label.
$ws.Range("A10").Value = "2020 (21.12.20)"

# Update the active selection to match the saved view state.
$ws.Range("K7").Select()
